$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B.
# This shifts the existing B,C,D,E (Jun_17,Jun_15,Jun_13,Jun_10) to E,F,G,H
# and keeps all styling/values of those columns (including the highlighted
# price-target cells) intact in their new position.
$ws.Range("B:D").Insert()

# Make sure all the date columns (C..H) keep the same custom width as
# before (8 characters).
for ($c = 3; $c -le 8; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 7.14285714
}

# New header values for the 3 newly inserted date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill in "UN" (unchanged rating) for the new columns on every existing
# analyst row (2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Add the new analyst group rows.
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"
$ws.Cells.Item(28, 4).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
$ws.Cells.Item(29, 4).Value = "UN"
